# "ajout 3 articles en txt"
#
# The "list articles" sheet had a few placeholder rows (8, 13, 14) that were
# completely empty in column B, plus row 17 whose title cell held only a
# stray space. This commit fills those in with three newly-read papers
# (Transfer Learning from Speaker Verification to Multispeaker TTS,
# FastSpeech, Deep Voice / Deep Voice 2) and marks them "ok" in the
# "manual rewriting" column, matching the formatting (wrapped text, slightly
# shorter row height) already used elsewhere on rows with long titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: was entirely empty -> new article title.
$ws.Cells.Item(8, 2).Value = "Transfer Learning from Speaker Verification to Multispeaker Text-To-Speech Synthesis"
$ws.Cells.Item(8, 2).WrapText = $true
$ws.Rows.Item(8).RowHeight = 15

# Row 13: was entirely empty -> new article title + review status.
$ws.Cells.Item(13, 2).Value = "FastSpeech: Fast, Robust and Controllable Text to Speech"
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Rows.Item(13).RowHeight = 15
$ws.Cells.Item(13, 3).Value = "ok"

# Row 14: was entirely empty -> new article title + review status.
$ws.Cells.Item(14, 2).Value = "Deep Voice: Real-time Neural Text-to-Speech"
$ws.Cells.Item(14, 2).WrapText = $true
$ws.Rows.Item(14).RowHeight = 15
$ws.Cells.Item(14, 3).Value = "ok"

# Row 17: title cell only had a stray space -> replace with real title + status.
$ws.Cells.Item(17, 2).Value = "Deep Voice 2: Multi-Speaker Neural Text-to-Speech"
$ws.Cells.Item(17, 2).WrapText = $true
$ws.Rows.Item(17).RowHeight = 15
$ws.Cells.Item(17, 3).Value = "ok"

# Row 19 (MELLOTRON) was missing its review status.
$ws.Cells.Item(19, 3).Value = "ok"

# Reflect where the author's cursor ended up after typing the new rows.
$ws.Range("C20").Select() | Out-Null
